# Updated symbol list on Sun Jan 29 22:22:12 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for the crypto rows,
# and fixes the BOLO / CoinbaseStockToken rows (48-49), which swapped
# positions and got new Price/Volume figures.
#
# Numeric-looking values are entered with a leading apostrophe so Excel
# keeps them as text (matching the sheet's existing inlineStr/text cells)
# instead of auto-converting them to Number/Percentage values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'318.52"
$ws.Range("E2").Value = "'4.18%"
$ws.Range("D3").Value = "'39.71"
$ws.Range("E3").Value = "'2.37%"
$ws.Range("D4").Value = "'5.139"
$ws.Range("E4").Value = "'0.50%"
$ws.Range("D5").Value = "'0.08215"
$ws.Range("E5").Value = "'1.86%"
$ws.Range("D6").Value = "'2.050"
$ws.Range("E6").Value = "'6.21%"
$ws.Range("D7").Value = "'8.373"
$ws.Range("E7").Value = "'4.47%"
$ws.Range("D8").Value = "'4.316"
$ws.Range("E8").Value = "'2.76%"
$ws.Range("D9").Value = "'0.9420"
$ws.Range("E9").Value = "'1.67%"
$ws.Range("D10").Value = "'0.1370"
$ws.Range("E10").Value = "'-4.98%"
$ws.Range("D11").Value = "'0.2005"
$ws.Range("E11").Value = "'4.97%"
$ws.Range("D12").Value = "'0.09088"
$ws.Range("E12").Value = "'1.23%"
$ws.Range("D13").Value = "'0.03521"
$ws.Range("E13").Value = "'0.37%"
$ws.Range("D14").Value = "'0.09806"
$ws.Range("E14").Value = "'0.30%"
$ws.Range("D15").Value = "'0.001412"
$ws.Range("E15").Value = "'0.98%"
$ws.Range("D16").Value = "'0.005971"
$ws.Range("E16").Value = "'0.41%"
$ws.Range("D17").Value = "'3.683"
$ws.Range("E17").Value = "'-1.93%"
$ws.Range("D18").Value = "'3.315"
$ws.Range("E18").Value = "'-0.43%"
$ws.Range("D19").Value = "'0.3482"
$ws.Range("E19").Value = "'0.59%"
$ws.Range("D20").Value = "'0.1322"
$ws.Range("E20").Value = "'-0.34%"
$ws.Range("D21").Value = "'4.953"
$ws.Range("E21").Value = "'5.97%"
$ws.Range("D22").Value = "'0.2449"
$ws.Range("E22").Value = "'1.45%"
$ws.Range("D23").Value = "'0.04365"
$ws.Range("E23").Value = "'-0.08%"
$ws.Range("D24").Value = "'0.001229"
$ws.Range("E24").Value = "'0.25%"
$ws.Range("D25").Value = "'0.004796"
$ws.Range("E25").Value = "'12.16%"
$ws.Range("E26").Value = "'-0.01%"
$ws.Range("D27").Value = "'0.0004001"
$ws.Range("E27").Value = "'-10.03%"
$ws.Range("D39").Value = "'0.02273"
$ws.Range("E39").Value = "'11.76%"
$ws.Range("D40").Value = "'0.05202"
$ws.Range("E40").Value = "'2.95%"
$ws.Range("D41").Value = "'0.007743"
$ws.Range("E41").Value = "'2.85%"
$ws.Range("D42").Value = "'0.009926"
$ws.Range("E42").Value = "'1.19%"
$ws.Range("D43").Value = "'0.1407"
$ws.Range("E43").Value = "'4.71%"
$ws.Range("D44").Value = "'0.002079"
$ws.Range("E44").Value = "'-0.66%"
$ws.Range("D45").Value = "'0.009660"
$ws.Range("E45").Value = "'-1.75%"
$ws.Range("D46").Value = "'0.00006612"
$ws.Range("E46").Value = "'6.47%"
$ws.Range("E47").Value = "'-0.12%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.001691"
$ws.Range("E48").Value = "'-6.26%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.002944"
$ws.Range("E49").Value = "'2.49%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.12%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.12%"
